$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = "Success"
$ws.Range("E3").Value = "08907097173605"
$ws.Range("E3").Style = "Normal"
$ws.Range("T3").Value = "Success"
$ws.Range("E4").Value = "08907097173612"
$ws.Range("E4").Style = "Normal"
$ws.Range("T4").Value = "Success"
$ws.Range("E5").Value = "08907097173629"
$ws.Range("E5").Style = "Normal"
$ws.Range("T5").Value = "Success"
$ws.Range("E6").Value = "08907097173636"
$ws.Range("E6").Style = "Normal"
$ws.Range("T6").Value = "Success"
